$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the A74 cell (hyperlinked URL) so it loses its HyperLink styling ---
# Remove the hyperlink object attached to A74.
$ws.Range("A74").Hyperlinks.Delete()

# Delete the whole row (drops the HyperLink-styled cell) and re-insert a blank
# row in its place; the new row inherits formatting ("source" style) from the
# row above it (A73, "Haut Commissariat Au Plan").
$ws.Rows.Item(74).Delete()
$ws.Rows.Item(74).Insert()
$ws.Range("A74").Value = "http://www.hcp.ma/pubData/RecensementEconomique/PublicationFinale/2000-2001.pdf"

# --- Insert blank separator rows between each "source" detail line ---
# Insert bottom-to-top so row numbers used below don't shift out from under us.
# Each new row inherits the style of the row immediately above it.
$ws.Rows.Item(75).Insert()
$ws.Rows.Item(74).Insert()
$ws.Rows.Item(73).Insert()

# --- Replace the long UN-ECA citation with a duplicate "UN-ECA" line ---
$ws.Range("A82").Value = "UN-ECA"
